# TC_52210_11_12_13 test data update
# - 40V_And_DC_Units (sheet2) and 40V_AC_DC_Units (sheet3):
#     Row 8's "Expected Power Calculation Text" changes from "NA" to an
#     actual "Load at 40V..." message (now wraps onto two lines), and the
#     comma-separated messages in rows 9-10 are reordered so the
#     "Load at 40V..." clause comes first instead of last.
# - 40V_AC_DC_Units becomes the active sheet/tab (was Extra_ISUnits).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 40V_And_DC_Units
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("40V_And_DC_Units")

$ws2.Range("D8").WrapText = $true
$ws2.Range("D8").Value = "'Load at 40V has reached 100% for panel Node1-Pro32xD-Built-in Loop-A"
$ws2.Rows.Item(8).RowHeight = 28.8

$ws2.Range("D9").Value = "'Load at 40V has reached 95% for panel Node1-Pro32xD-Built-in Loop-A,DC value has reached 95% for panel Node1-Pro32xD,"

$ws2.Range("D10").Value = "'Load at 40V has reached 100% for panel Node1-Pro32xD-Built-in Loop-A,DC value has reached 100% for panel Node1-Pro32xD"

# Selection moves to D9 (sheet stays inactive once sheet 3 is selected below)
$ws2.Range("D9").Select()

# ---------------------------------------------------------------------
# 40V_AC_DC_Units
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("40V_AC_DC_Units")

$ws3.Range("D8").WrapText = $true
$ws3.Range("D8").Value = "'Load at 40V has reached 100% for panel Node1-Pro32xD-Built-in Loop-A"
$ws3.Rows.Item(8).RowHeight = 28.8

$ws3.Range("D9").Value = "'Load at 40V has reached 100% for panel Node1-Pro32xD-Built-in Loop-A,DC value has reached 100% for panel Node1-Pro32xD,AC value has reached 95% for panel Node1-Pro32xD"

$ws3.Range("D10").Value = "'Load at 40V has reached 100% for panel Node1-Pro32xD-Built-in Loop-A,DC value has reached 100% for panel Node1-Pro32xD,AC value has reached 100% for panel Node1-Pro32xD"

# Selecting D10 here both sets the new selection AND makes this sheet the
# active tab (matching workbookView activeTab moving from Extra_ISUnits to
# this sheet, and tabSelected moving off Extra_ISUnits onto this sheet).
$ws3.Range("D10").Select()
